$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date.
$ws.Name = "Through 2022-10-20"

# Update the October row header label.
$ws.Range("A11").Value = "October (through 10-20)"

# Update October (row 11) values.
$ws.Range("B11").Value = 19
$ws.Range("C11").Value = 32
$ws.Range("D11").Value = 35
$ws.Range("E11").Value = 48
$ws.Range("F11").Value = 31
$ws.Range("G11").Value = 93
$ws.Range("H11").Value = 128
$ws.Range("I11").Value = 67

# Update Total (row 12) values.
$ws.Range("B12").Value = 245
$ws.Range("C12").Value = 461
$ws.Range("D12").Value = 662
$ws.Range("E12").Value = 596
$ws.Range("F12").Value = 453
$ws.Range("G12").Value = 994
$ws.Range("H12").Value = 1375
$ws.Range("I12").Value = 1344
